$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '66.041.69'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -1.25%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.750.83'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +0.88%  '
$ws.Range("E4").Value = '  +0.15%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '403.99'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -4.71%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '129.00'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -2.36%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '3.738.65'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.88%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.602'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -6.17%  '
$ws.Range("E9").Value = '  -0.06%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.717'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -6.66%  '
$ws.Range("E11").Value = '  -9.45%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0000356'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -7.57%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '40.30'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -6.13%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '4.334.58'
$ws.Range("D14").Style = "Normal"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '9.60'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -6.34%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '14.42'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +9.60%  '
$ws.Range("E17").Value = '  -1.81%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '3.741.81'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.12%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '19.34'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -7.65%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '66.391.60'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.75%  '
$ws.Range("E21").Value = '  -7.10%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '407.27'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -9.26%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '14.35'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -9.83%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '84.82'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -6.49%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.99'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -6.33%  '
$ws.Range("B26").Value = 'LEO'
$ws.Range("C26").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '5.64'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +13.55%  '
$ws.Range("B27").Value = 'EthereumClassic'
$ws.Range("C27").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '35.86'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -6.23%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '3.08'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -6.97%  '
$ws.Range("E29").Value = '  -9.53%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '12.31'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -2.59%  '
$ws.Range("E31").Value = '  -2.95%  '
$ws.Range("E32").Value = '  -3.96%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '7.09'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -2.79%  '
$ws.Range("E34").Value = '  -5.62%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '38.73'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -8.57%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.00'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +0.00%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '54.97'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -2.60%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.0₃0720'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -2.39%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0452'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -7.97%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.87'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -6.87%  '
$ws.Range("E41").Value = '  +0.14%  '
$ws.Range("B42").Value = 'EnergySwap'
$ws.Range("C42").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '27.11'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -6.57%  '
$ws.Range("B43").Value = 'Stellar'
$ws.Range("C43").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.134'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -9.08%  '
$ws.Range("B44").Value = 'ApeXProtocol'
$ws.Range("C44").Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '3.18'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +21.09%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '144.90'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.76%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.22'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -7.22%  '
$ws.Range("E47").Value = '  -4.51%  '
$ws.Range("E48").Value = '  -4.03%  '
$ws.Range("E49").Value = '  -4.46%  '
$ws.Range("E50").Value = '  -4.77%  '
$ws.Range("E51").Value = '  -6.54%  '
